$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.142.93"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.24%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.825.23"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "312.55"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4635"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -1.61%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3630"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07297"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.69%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8703"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.37%  "

$ws.Range("E11").Value = "  -1.68%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.875.35"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.07641"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.25%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.342"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.49%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "92.41"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.469"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.009"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008629"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "27.445.34"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.49"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "

$ws.Range("E22").Value = "  -1.82%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.56"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.26%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.094.07"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.92%  "

$ws.Range("E25").Value = "  -1.15%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "151.13"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.23"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.14%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.086"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.71%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.107"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.31%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "116.19"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.963"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.7367"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.06%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.145"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.38%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.452"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.16%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.012"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.502"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.25%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.081"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05243"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01915"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.44%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.928"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -2.75%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "7.158"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.5200"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.92%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.1626"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.29%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.296"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.89%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.4839"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.35%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.011"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.35%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "103.52"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.636"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06270"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
